$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 1836
$ws.Range("I3").Value = 1943
$ws.Range("D4").Value = 1924
$ws.Range("H4").Value = 1659
$ws.Range("I4").Value = 493
$ws.Range("I5").Value = 169
$ws.Range("I6").Value = 2378
$ws.Range("D7").Value = 28114
$ws.Range("H7").Value = 25971
$ws.Range("I7").Value = 6819

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I4").Value = 5
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 91

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 65
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 216

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 42
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 63
$ws.Range("I3").Value = 91
$ws.Range("I7").Value = 265

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 41
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 68
$ws.Range("I7").Value = 238
$ws.Range("I8").Value = 421
$ws.Range("I10").Value = 53
$ws.Range("I11").Value = 120
$ws.Range("I13").Value = 8
$ws.Range("I14").Value = 36
$ws.Range("I18").Value = 54
$ws.Range("I19").Value = 199
$ws.Range("I20").Value = 187
$ws.Range("I29").Value = 437
$ws.Range("I32").Value = 12
$ws.Range("I33").Value = 318
$ws.Range("I35").Value = 8
$ws.Range("I36").Value = 88
$ws.Range("I37").Value = 216
$ws.Range("I40").Value = 10
$ws.Range("I42").Value = 226
$ws.Range("I43").Value = 61
$ws.Range("I47").Value = 52
$ws.Range("I48").Value = 67
$ws.Range("I49").Value = 40
$ws.Range("I51").Value = 67
$ws.Range("I52").Value = 140
$ws.Range("I53").Value = 71
$ws.Range("I54").Value = 149
$ws.Range("I55").Value = 78
$ws.Range("I58").Value = 4
$ws.Range("D63").Value = 316
$ws.Range("H63").Value = 187
$ws.Range("I63").Value = 34
$ws.Range("I65").Value = 159
$ws.Range("I67").Value = 265
$ws.Range("I70").Value = 12
$ws.Range("I71").Value = 14
$ws.Range("I73").Value = 64
$ws.Range("I76").Value = 111
$ws.Range("I78").Value = 91
$ws.Range("I80").Value = 22
$ws.Range("I83").Value = 129
$ws.Range("I85").Value = 324
$ws.Range("I86").Value = 42
$ws.Range("I90").Value = 81
$ws.Range("I95").Value = 114
$ws.Range("I96").Value = 91
$ws.Range("I97").Value = 56
$ws.Range("I99").Value = 122
$ws.Range("D101").Value = 28114
$ws.Range("H101").Value = 25971
$ws.Range("I101").Value = 6819

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 44
$ws.Range("I3").Value = 51
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 76
$ws.Range("I3").Value = 109
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 318

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 131
$ws.Range("I3").Value = 140
$ws.Range("I4").Value = 14
$ws.Range("I5").Value = 12
$ws.Range("I6").Value = 140
$ws.Range("I7").Value = 437

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 49
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 199

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 6
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 20
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 111

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 78
$ws.Range("I3").Value = 126
$ws.Range("I6").Value = 93
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 77
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("I4").Value = 2
$ws.Range("I6").Value = 8

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 18
$ws.Range("I3").Value = 23
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 91

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 45
$ws.Range("I3").Value = 56
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 187

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 11
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 54

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 24
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 57
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 15
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 22
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 51
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("I2").Value = 3
$ws.Range("I7").Value = 8

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 22
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 12

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 12

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 125
$ws.Range("I3").Value = 114
$ws.Range("I4").Value = 25
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 421

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 28
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 26
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 81

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 14

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I3").Value = 3
$ws.Range("I7").Value = 10

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 70
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 238

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("I6").Value = 2
$ws.Range("I7").Value = 4
